$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.816.01"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "1.633.71"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5037"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06419"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07701"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.251"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.627.52"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "1.858.93"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5461"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "25.829.06"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.46%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.951"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.969"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.916"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1148"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.712"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05029"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.187"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.536"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.358"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "1.177.63"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8955"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.604"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5609"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.547"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.667"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8068"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "1.771.28"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("E51").Value = "  -0.42%  "
